$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("A15").Value = 45975

$ws.Range("B15").Value = 64

$ws.Range("A15:B15").Select()
